# RP3 FLT_EFF 2020 workbook update — "202106 release uk" (#41)
#   * UK removed in RP3 (SES area no longer includes the UK)
#   * download files updated (refreshed figures for the June 2021 release)
#   * Updates June release with May data
#
# This script reproduces the semantic edits described by the commit:
#  1. FLT_EFF_YY  : bump release date, refresh yearly SES-area figures
#  2. FLT_EFF_MM  : bump release date, rename "SES AREA" -> "SES AREA (RP3)",
#                   refresh monthly SES-area figures (Jan-19..Dec-20)
#  3. ERT_FLT_EFF_LOC : drop the "UK (Continental)" row (UK removed from SES area)
#  4. Change Log  : add a log entry documenting the UK removal

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) FLT_EFF_YY
# ---------------------------------------------------------------------------
$wsYY = $wb.Worksheets.Item("FLT_EFF_YY")

$wsYY.Range("B2").Value = 44351

$wsYY.Range("D6").Value = 0.0474
$wsYY.Range("F6").Value = 0.0284

$wsYY.Range("D7").Value = 0.0457
$wsYY.Range("F7").Value = 0.0268

$wsYY.Range("D8").Value = 0.0456
$wsYY.Range("F8").Value = 0.0272

$wsYY.Range("D9").Value = 0.0453
$wsYY.Range("F9").Value = 0.0284

$wsYY.Range("D10").Value = 0.0438
$wsYY.Range("E10").Value = 0.0398
$wsYY.Range("F10").Value = 0.0251

# ---------------------------------------------------------------------------
# 2) FLT_EFF_MM
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("FLT_EFF_MM")

$wsMM.Range("B2").Value = 44351

# SES area was redefined for RP3 (UK removed) -> relabel the header
$wsMM.Range("B4").Value = "SES AREA (RP3)"

# Refreshed monthly figures, Jan-19 .. Dec-20 (rows 6-29)
$monthly = @(
    @{ Row = 6;  B = 0.0456; D = 0.0273 }
    @{ Row = 7;  B = 0.0456; D = 0.0274 }
    @{ Row = 8;  B = 0.0455; D = 0.0275 }
    @{ Row = 9;  B = 0.0454; D = 0.0275 }
    @{ Row = 10; B = 0.0454; D = 0.0275 }
    @{ Row = 11; B = 0.0453; D = 0.0276 }
    @{ Row = 12; B = 0.0453; D = 0.0278 }
    @{ Row = 13; B = 0.0453; D = 0.028  }
    @{ Row = 14; B = 0.0453; D = 0.0281 }
    @{ Row = 15; B = 0.0453; D = 0.0283 }
    @{ Row = 16; B = 0.0452; D = 0.0283 }
    @{ Row = 17; B = 0.0453; D = 0.0284 }
    @{ Row = 18; B = 0.0453; C = 0.0418; D = 0.0285 }
    @{ Row = 19; B = 0.0452; C = 0.0417; D = 0.0285 }
    @{ Row = 20; B = 0.0451; C = 0.0416; D = 0.0283 }
    @{ Row = 21; B = 0.0451; C = 0.0415; D = 0.0283 }
    @{ Row = 22; B = 0.0449; C = 0.0413; D = 0.0281 }
    @{ Row = 23; B = 0.0448; C = 0.0411; D = 0.0277 }
    @{ Row = 24; B = 0.0446; C = 0.0409; D = 0.0272 }
    @{ Row = 25; B = 0.0444; C = 0.0407; D = 0.0268 }
    @{ Row = 26; B = 0.0442; C = 0.0404; D = 0.0263 }
    @{ Row = 27; B = 0.044;  C = 0.0401; D = 0.0258 }
    @{ Row = 28; B = 0.044;  C = 0.04;   D = 0.0255 }
    @{ Row = 29; B = 0.0438; C = 0.0398; D = 0.0251 }
)

foreach ($m in $monthly) {
    $wsMM.Range("B" + $m.Row).Value = $m.B
    if ($m.ContainsKey("C")) {
        $wsMM.Range("C" + $m.Row).Value = $m.C
    }
    $wsMM.Range("D" + $m.Row).Value = $m.D
}

# ---------------------------------------------------------------------------
# 3) ERT_FLT_EFF_LOC — UK (Continental) removed from the SES area entirely
# ---------------------------------------------------------------------------
$wsLOC = $wb.Worksheets.Item("ERT_FLT_EFF_LOC")

$wsLOC.Rows.Item(34).Delete()

# Conditional-formatting ranges tracked F6:F34; shrink them to match.
$fcs = $wsLOC.Range("F6:F33").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($wsLOC.Range("F6:F33"))
}

# ---------------------------------------------------------------------------
# 4) Change Log — record the UK removal
# ---------------------------------------------------------------------------
$wsLog = $wb.Worksheets.Item("Change Log")

$wsLog.Range("A2").Value = 44351
$wsLog.Range("B2").Value = "UK"
$wsLog.Range("D2").Value = "UK removed from SES area"
